$d = $word.ActiveDocument

# Keep Word from "smartening" straight quotes/apostrophes into curly ones
# while we run replacements below.
$word.Options.AutoFormatAsYouTypeReplaceQuotes = $false
$word.Options.AutoFormatReplaceQuotes = $false

function Replace-Text($findText, $replaceText) {
    $rng = $d.Content
    $ok = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $ok) {
        throw ("Replace-Text: text not found: " + $findText)
    }
}

# Insert $newText right after the first occurrence of $afterText (plain insert;
# it will likely get merged visually into the neighbouring same-format run --
# callers that need a hard run boundary should re-split with Split-Runs afterwards).
function Insert-After($afterText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($afterText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw ("Insert-After: text not found: " + $afterText)
    }
    $rng.Collapse(0) | Out-Null
    $rng.InsertBefore($newText)
}

# Force paragraph $para's text to be split into runs at the boundaries implied by
# the given ordered list of $segments (whose concatenation must equal the
# paragraph's text, sans the trailing paragraph mark). Toggling a character
# property on a precise sub-range and reverting it is enough to make the engine
# keep that sub-range as its own <w:r>, instead of silently re-coalescing it into
# its same-formatted neighbours.
function Split-Runs($para, [string[]]$segments) {
    $pStart = $para.Range.Start
    $offset = 0
    foreach ($seg in $segments) {
        $len = $seg.Length
        if ($len -gt 0) {
            $r = $d.Range($pStart + $offset, $pStart + $offset + $len)
            $r.Font.Bold = $true
            $r.Font.Bold = $false
        }
        $offset += $len
    }
}

# ---------------------------------------------------------------------------
# 1. Font fix: TimesNewToman -> Times New Roman, everywhere.
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    $p.Range.Font.Name = "Times New Roman"
}

# ---------------------------------------------------------------------------
# 2. Title
# ---------------------------------------------------------------------------
Replace-Text "Unraveling the Enigma of Black Holes" "Biology: Exploring the Symphony of Life"

# ---------------------------------------------------------------------------
# 3. Author name
# ---------------------------------------------------------------------------
Replace-Text " Amelia Carter" " Emily Jones"

# ---------------------------------------------------------------------------
# 4. Email address: split "ameliacarter@cosmosinstitute" + "." + "org"
#    into "emily" + "." + "jones@eduworld" + "." + "org"
# ---------------------------------------------------------------------------
Replace-Text "ameliacarter@cosmosinstitute" "emily"
Insert-After "emily." "jones@eduworld."

$emailPara = $d.Paragraphs(3)
Split-Runs $emailPara @("emily", ".", "jones@eduworld", ".", "org")

# ---------------------------------------------------------------------------
# 5. Body paragraph
# ---------------------------------------------------------------------------
Replace-Text "In the boundless expanse of the universe, there lies a celestial enigma that has captivated the minds of scientists and lay people alike: black holes" "Biology, the study of life, is a journey that unravels the mysteries of the living world"

Replace-Text " These enigmatic entities are cosmic vacuums with an infinitely strong gravitational pull, from which nothing, not even light, can escape" " It's an exploration into the intricate workings of organisms, from the smallest microbes to the largest whales, and the dynamic interactions between them"

Replace-Text " They represent a dark frontier in our understanding of space and time, inviting us to explore their extraordinary properties and unravel their cosmic mysteries" " We embark on a quest to comprehend the secrets of life, delving into the symphony of biological processes that orchestrate the harmony of existence"

Insert-After "harmony of existence." " As we dissect the molecular mechanisms that govern cellular functions, we unravel the mysteries of genetics, the blueprint of life."
Insert-After "blueprint of life." " We investigate the intricate web of ecosystems, the intricate relationships between organisms and their environments, and the delicate balance that sustains our planet's biodiversity."

Replace-Text "In this exploration, we embark on a journey to penetrate the veil of darkness surrounding black holes" "Biology propels us into the realm of human biology, unraveling the complexities of our bodies, from the microscopic world of cells to the sophisticated systems that regulate our thoughts, emotions, and actions"

Replace-Text " We begin by unraveling their gravitational peculiarities, delving into the concept of event horizons and their role in creating an inescapable boundary" " We delve into the wonders of the immune system, our body's intricate defense mechanism, fending off invaders and maintaining internal harmony"

Replace-Text " Additionally, we examine the mind-bending phenomena occurring near black holes, including time dilation and the mesmerizing behavior of light" " We explore the marvels of reproduction, the miracle of new life emerging from the fusion of genetic material"

Insert-After "fusion of genetic material." " By studying the human body, we gain a profound appreciation for the resilience, adaptability, and interconnectedness of life."

Replace-Text "Furthermore, we investigate the birth of these cosmic leviathans through the death of massive stars" "Furthermore, biology illuminates the interconnectedness of life on Earth"

Replace-Text " We analyze the various evolutionary pathways, such as stellar collapse and supernova explosions, that lead to the formation of these enigmatic entities" " We uncover the intricate web of interdependence among organisms, the delicate balance of ecosystems, and the profound impact of human activities on the natural world"

Replace-Text " Our quest for knowledge leads us to question the ultimate fate of black holes, considering their hypothetical evaporation through Hawking radiation and their possible involvement in mysterious cosmic phenomena like gravitational waves" " Biology empowers us to recognize our responsibility as stewards of the planet, inspiring us to act as conscientious citizens, preserving and protecting the diversity of life for generations to come"

$VT = [char]11
$bodyPara = $d.Paragraphs(5)
Split-Runs $bodyPara @(
    "Biology, the study of life, is a journey that unravels the mysteries of the living world",
    ".",
    " It's an exploration into the intricate workings of organisms, from the smallest microbes to the largest whales, and the dynamic interactions between them",
    ".",
    " We embark on a quest to comprehend the secrets of life, delving into the symphony of biological processes that orchestrate the harmony of existence",
    ".",
    " As we dissect the molecular mechanisms that govern cellular functions, we unravel the mysteries of genetics, the blueprint of life",
    ".",
    " We investigate the intricate web of ecosystems, the intricate relationships between organisms and their environments, and the delicate balance that sustains our planet's biodiversity",
    ".",
    "$VT",
    "$VT" + "Biology propels us into the realm of human biology, unraveling the complexities of our bodies, from the microscopic world of cells to the sophisticated systems that regulate our thoughts, emotions, and actions",
    ".",
    " We delve into the wonders of the immune system, our body's intricate defense mechanism, fending off invaders and maintaining internal harmony",
    ".",
    " We explore the marvels of reproduction, the miracle of new life emerging from the fusion of genetic material",
    ".",
    " By studying the human body, we gain a profound appreciation for the resilience, adaptability, and interconnectedness of life",
    ".",
    "$VT",
    "$VT" + "Furthermore, biology illuminates the interconnectedness of life on Earth",
    ".",
    " We uncover the intricate web of interdependence among organisms, the delicate balance of ecosystems, and the profound impact of human activities on the natural world",
    ".",
    " Biology empowers us to recognize our responsibility as stewards of the planet, inspiring us to act as conscientious citizens, preserving and protecting the diversity of life for generations to come",
    "."
)

# ---------------------------------------------------------------------------
# 6. Summary paragraph
# ---------------------------------------------------------------------------
Replace-Text "Black holes stand as a testament to the vastness and complexity of the universe, beckoning us to push the boundaries of our scientific understanding" "Biology is a captivating subject that delves into the intricacies of life, unraveling the mysteries of living organisms and their interactions with each other and their environment"

Replace-Text " Through an examination of their gravitational anomalies, formation mechanisms, and cosmic interactions, we have delved into the depths of these celestial conundrums" " By exploring the symphony of biological processes, we gain a deeper understanding of our bodies, appreciate the marvels of life's diversity, and recognize our role as stewards of the planet"

Replace-Text " As we continue to unravel the enigma of black holes, we unlock new insights into the fundamental nature of space, time, and the dynamics of the cosmos" " Biology inspires us to question, investigate, and marvel at the wonders of the natural world, nurturing a lifelong appreciation for the beauty and complexity of life"

$summaryPara = $d.Paragraphs(7)
Split-Runs $summaryPara @(
    "Biology is a captivating subject that delves into the intricacies of life, unraveling the mysteries of living organisms and their interactions with each other and their environment",
    ".",
    " By exploring the symphony of biological processes, we gain a deeper understanding of our bodies, appreciate the marvels of life's diversity, and recognize our role as stewards of the planet",
    ".",
    " Biology inspires us to question, investigate, and marvel at the wonders of the natural world, nurturing a lifelong appreciation for the beauty and complexity of life",
    "."
)

# ---------------------------------------------------------------------------
# 7. New trailing empty paragraph
# ---------------------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0) | Out-Null
$endRange.InsertParagraphAfter()

Write-Output "done"
